$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ID_db_Lemma (column C) values that shifted due to corrected lookup generation
$ws.Range("C44").Value = 1000043
$ws.Range("C45").Value = 1000044
$ws.Range("C46").Value = 1000045
$ws.Range("C47").Value = 1000046
$ws.Range("C48").Value = 1000047
$ws.Range("C49").Value = 1000049
$ws.Range("C50").Value = 1000050
$ws.Range("C51").Value = 1000051
$ws.Range("C52").Value = 1000052
$ws.Range("C53").Value = 1000053
$ws.Range("C54").Value = 1000054
$ws.Range("C55").Value = 1000055
$ws.Range("C56").Value = 1000056
$ws.Range("C57").Value = 1000057
$ws.Range("C58").Value = 1000058
$ws.Range("C59").Value = 1000059
$ws.Range("C60").Value = 1000060
$ws.Range("C61").Value = 1000061
$ws.Range("C62").Value = 1000062
$ws.Range("C63").Value = 1000063
$ws.Range("C64").Value = 1000064
$ws.Range("C65").Value = 1000065
$ws.Range("C66").Value = 1000066
$ws.Range("C67").Value = 1000068
$ws.Range("C68").Value = 1000069
$ws.Range("C69").Value = 1000071
$ws.Range("C70").Value = 1000073
$ws.Range("C71").Value = 1000074

# Update Lemma (column E) text values that were re-mapped/corrected
$ws.Range("E27").Value = "Centro operativo regionale"
$ws.Range("E35").Value = "SOREU dei laghi"
$ws.Range("E36").Value = "Vigili del fuoco"
$ws.Range("E37").Value = "Nucleo Unitario di Valutazione e Risposta Emergenze transfrontaliere"
$ws.Range("E38").Value = "Sala Operativa Regionale dell'Emergenza Urgenza"
$ws.Range("E39").Value = "sezione del militare e della protezione della popolazione"
$ws.Range("E40").Value = "Centrale nazionale d'allarme"
$ws.Range("E41").Value = "Stato maggiore federale Protezione della popolazione"
$ws.Range("E42").Value = "Protezione civile"
$ws.Range("E43").Value = "Legge federale sulla protezione della popolazione e sulla protezione civile"
$ws.Range("E49").Value = "Struttura operativa"
$ws.Range("E50").Value = "Organizzazione partner"
$ws.Range("E52").Value = "Emergenza / Stato di emergenza / Evento emergenziale / Evento"
$ws.Range("E53").Value = "Stato di necessità"
$ws.Range("E68").Value = "Legge sulla protezione della popolazione del 26 febbraio 2007"

# Remove the two trailing rows that no longer exist after the fix (old rows 72 and 73)
$ws.Rows.Item(73).Delete()
$ws.Rows.Item(72).Delete()

Write-Output "edit complete"
